# Rows 32-39 of the "Artfynd" sheet had their entire data-row contents
# permuted (observations re-ordered / re-matched to different coordinate
# rows), while the row numbers themselves stayed put. Row 38 is unchanged.
#
# Mapping: destination row -> source row (i.e. "row X should end up
# containing what row Y currently contains"):
#   32 <- 34
#   33 <- 36
#   34 <- 35
#   35 <- 37
#   36 <- 33
#   37 <- 39
#   38 <- 38   (unchanged)
#   39 <- 32

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every column (within A:AY) that is populated on at least one of rows 32-39.
$cols = @('A','B','C','D','E','F','G','H','I','K','P','Q','R','S','T','U','V','W', `
          'Y','Z','AA','AB','AC','AD','AE','AG','AR','AT','AW','AX','AY')

$sourceRows = 32..39

# 1) Snapshot every used cell of every source row BEFORE any writes happen,
#    since destinations and sources overlap (same 8 rows, permuted).
$snapshot = @{}
foreach ($r in $sourceRows) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $addr = "$col$r"
        $cell = $ws.Range($addr)
        $rowVals[$col] = $cell.Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Destination -> source row mapping.
$mapping = @{
    32 = 34
    33 = 36
    34 = 35
    35 = 37
    36 = 33
    37 = 39
    38 = 38
    39 = 32
}

# 3) Which columns exist (are non-blank) per *source* row, so we know which
#    columns must be cleared on the destination if the incoming row doesn't
#    use them (e.g. K/AC/AR only appear on some of these rows).
$presentCols = @{
    32 = @('A','B','C','D','E','F','G','H','I','P','Q','R','S','T','U','V','W','Y','Z','AA','AB','AD','AE','AG','AT','AW','AX','AY')
    33 = @('A','B','C','D','E','F','G','H','I','P','Q','R','S','T','U','V','W','Y','Z','AA','AB','AD','AE','AG','AT','AW','AX','AY')
    34 = @('A','B','C','D','E','F','G','H','I','P','Q','R','S','T','U','V','W','Y','Z','AA','AB','AD','AE','AG','AT','AW','AX','AY')
    35 = @('A','B','C','D','E','F','G','H','I','K','P','Q','R','S','T','U','V','W','Y','Z','AA','AB','AD','AE','AG','AR','AT','AW','AX','AY')
    36 = @('A','B','C','D','E','F','G','H','I','P','Q','R','S','T','U','V','W','Y','Z','AA','AB','AD','AE','AG','AT','AW','AX','AY')
    37 = @('A','B','C','D','E','F','G','H','I','K','P','Q','R','S','T','U','V','W','Y','Z','AA','AB','AC','AD','AE','AG','AT','AW','AX','AY')
    38 = @('A','B','C','D','E','F','G','H','I','P','Q','R','S','T','U','V','W','Y','Z','AA','AB','AD','AE','AG','AT','AW','AX','AY')
    39 = @('A','B','C','D','E','F','G','H','I','P','Q','R','S','T','U','V','W','Y','Z','AA','AB','AD','AE','AG','AT','AW','AX','AY')
}

# 4) Write each destination row from the matching snapshot, clearing any
#    column that the incoming (source) row doesn't actually populate.
foreach ($destRow in $sourceRows) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    $srcPresent = $presentCols[$srcRow]

    foreach ($col in $cols) {
        $addr = "$col$destRow"
        if ($srcPresent -contains $col) {
            $ws.Range($addr).Value = $srcVals[$col]
        } else {
            $ws.Range($addr).ClearContents()
        }
    }
}
